$wb = $excel.ActiveWorkbook

# --- Sheets: rename the existing sheet, add a new one after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "compradores"

$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws2.Name = "apidata"

# --- Sheet "compradores" data ---
$ws1.Range("A1").Value = "dataKey"
$ws1.Range("B1").Value = "nombre"
$ws1.Range("C1").Value = "apellido"
$ws1.Range("D1").Value = "email"
$ws1.Range("E1").Value = "genero"

$ws1.Range("A2").Value = "comprador1"
$ws1.Range("B2").Value = "Giulio"
$ws1.Range("C2").Value = "Faragalli"
$ws1.Range("D2").Value = "GFemail@email.com"
$ws1.Range("E2").Value = "masculino"

$ws1.Hyperlinks.Add($ws1.Range("D2"), "mailto:GFemail@email.com")

$ws1.Columns.Item(1).ColumnWidth = 12.1

$ws1.PageSetup.Orientation = 1

# --- Sheet "apidata" data ---
$ws2.Range("A1").Value = "dataKey"
$ws2.Range("B1").Value = "name"
$ws2.Range("C1").Value = "language"
$ws2.Range("D1").Value = "address"

$ws2.Range("A2").Value = "place1"
$ws2.Range("B2").Value = "AAhouse"
$ws2.Range("C2").Value = "English"
$ws2.Range("D2").Value = "World cross center"

$ws2.Range("A3").Value = "place2"
$ws2.Range("B3").Value = "BBhouse"
$ws2.Range("C3").Value = "Spanish"
$ws2.Range("D3").Value = "Sea cross center"

$ws2.PageSetup.Orientation = 1

# --- Selections: match the saved view state (compradores shows A1:E2,
#     apidata is the active sheet with E11 selected) ---
[void]$ws1.Range("A1:E2").Select()
[void]$ws2.Range("E11").Select()
[void]$ws2.Activate()
